# Restated IFRS-consolidated figures (values re-expressed, e.g. EPS/BPS/etc. corrected).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update restated financial figures for rows 2-6 (FY2014-FY2018 IFRS-consolidated) ---
# Row 2
$ws.Range("D2").Value = 10128
$ws.Range("E2").Value = -446
$ws.Range("F2").Value = -446
$ws.Range("G2").Value = -1281
$ws.Range("H2").Value = -1258
$ws.Range("I2").Value = -1258
$ws.Range("K2").Value = 14864
$ws.Range("L2").Value = 14887
$ws.Range("M2").Value = -22
$ws.Range("N2").Value = -22
$ws.Range("P2").Value = 2858
$ws.Range("Q2").Value = -827
$ws.Range("R2").Value = 1319
$ws.Range("S2").Value = -327
$ws.Range("T2").Value = 76
$ws.Range("U2").Value = -903
$ws.Range("V2").Value = 9049
$ws.Range("W2").Value = -4.4
$ws.Range("X2").Value = -12.42
$ws.Range("Y2").Value = -207.5
$ws.Range("Z2").Value = -7.61
$ws.Range("AA2").Value = -66595.69
$ws.Range("AB2").Value = -111.18
$ws.Range("AC2").Value = -1411556
$ws.Range("AD2").Value = -1.28
$ws.Range("AE2").Value = -25099
$ws.Range("AF2").Value = -72.02
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 89116

# Row 3
$ws.Range("D3").Value = 10572
$ws.Range("E3").Value = 78
$ws.Range("F3").Value = 78
$ws.Range("G3").Value = 109
$ws.Range("H3").Value = -49
$ws.Range("I3").Value = -49
$ws.Range("K3").Value = 14198
$ws.Range("L3").Value = 13131
$ws.Range("M3").Value = 1067
$ws.Range("N3").Value = 1067
$ws.Range("P3").Value = 649
$ws.Range("Q3").Value = 39
$ws.Range("R3").Value = 23
$ws.Range("S3").Value = 379
$ws.Range("T3").Value = 94
$ws.Range("U3").Value = -55
$ws.Range("V3").Value = 8469
$ws.Range("W3").Value = 0.74
$ws.Range("X3").Value = -0.47
$ws.Range("Y3").Value = -9.41
$ws.Range("Z3").Value = -0.34
$ws.Range("AA3").Value = 1230.19
$ws.Range("AB3").Value = -32.13
$ws.Range("AC3").Value = -18573
$ws.Range("AD3").Value = -48.76
$ws.Range("AE3").Value = 332022
$ws.Range("AF3").Value = 2.73
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 139392

# Row 4
$ws.Range("D4").Value = 4165
$ws.Range("E4").Value = -2832
$ws.Range("F4").Value = -2832
$ws.Range("G4").Value = -6975
$ws.Range("H4").Value = -7183
$ws.Range("I4").Value = -7183
$ws.Range("K4").Value = 8499
$ws.Range("L4").Value = 14613
$ws.Range("M4").Value = -6114
$ws.Range("N4").Value = -6114
$ws.Range("P4").Value = 649
$ws.Range("Q4").Value = -69
$ws.Range("R4").Value = 250
$ws.Range("S4").Value = -263
$ws.Range("T4").Value = 40
$ws.Range("U4").Value = -109
$ws.Range("V4").Value = 7911
$ws.Range("W4").Value = -68
$ws.Range("X4").Value = -172.47
$ws.Range("Y4").Value = 284.67
$ws.Range("Z4").Value = -63.29
$ws.Range("AA4").Value = -239.01
$ws.Range("AB4").Value = -1138.05
$ws.Range("AC4").Value = -2386566
$ws.Range("AD4").Value = -0.08
$ws.Range("AE4").Value = -2376350
$ws.Range("AF4").Value = -0.08
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 148554

# Row 5
$ws.Range("D5").Value = 1852
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 7035
$ws.Range("H5").Value = 7355
$ws.Range("I5").Value = 7355
$ws.Range("K5").Value = 7443
$ws.Range("L5").Value = 5988
$ws.Range("M5").Value = 1456
$ws.Range("N5").Value = 1456
$ws.Range("P5").Value = 617
$ws.Range("Q5").Value = 74
$ws.Range("R5").Value = -389
$ws.Range("S5").Value = -76
$ws.Range("T5").Value = 17
$ws.Range("U5").Value = 57
$ws.Range("V5").Value = 1230
$ws.Range("W5").Value = 0.24
$ws.Range("X5").Value = 397.04
$ws.Range("Y5").Value = -315.77
$ws.Range("Z5").Value = 92.27
$ws.Range("AA5").Value = 411.37
$ws.Range("AB5").Value = 524.4299999999999
$ws.Range("AC5").Value = 261718
$ws.Range("AD5").Value = 0.14
$ws.Range("AE5").Value = 47240
$ws.Range("AF5").Value = 0.78
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 3081168

# Row 6
$ws.Range("D6").Value = 2148
$ws.Range("E6").Value = -145
$ws.Range("F6").Value = -145
$ws.Range("G6").Value = -1438
$ws.Range("H6").Value = -1047
$ws.Range("I6").Value = -1047
$ws.Range("K6").Value = 5672
$ws.Range("L6").Value = 4133
$ws.Range("M6").Value = 1539
$ws.Range("N6").Value = 1539
$ws.Range("P6").Value = 702
$ws.Range("Q6").Value = 116
$ws.Range("R6").Value = 692
$ws.Range("S6").Value = -353
$ws.Range("T6").Value = 13
$ws.Range("U6").Value = 103
$ws.Range("V6").Value = 2164
$ws.Range("W6").Value = -6.75
$ws.Range("X6").Value = -48.75
$ws.Range("Y6").Value = -69.95
$ws.Range("Z6").Value = -15.97
$ws.Range("AA6").Value = 268.62
$ws.Range("AB6").Value = 513.61
$ws.Range("AC6").Value = -17662
$ws.Range("AD6").Value = -0.32
$ws.Range("AE6").Value = 5479
$ws.Range("AF6").Value = 1.03
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 28081731

# --- 당기순이익(비지배) (J) and 자본총계(비지배) (O) are no longer reported for rows 2-5 ---
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# --- Drop the unreliable (E) estimate rows 7-9 figures, keep only the year labels ---
$ws.Range("D7:AI9").ClearContents()
